$wb = $excel.ActiveWorkbook

# --- About sheet: add "Oregon" label next to the title in B1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B1").Value = "Oregon"

# --- Cost Data sheet: fix the /10 scaling bug on the annual totals ---
$wsCost = $wb.Worksheets.Item("Cost Data")

# Match the number format already used by neighboring "$"#,##0 cells
# (e.g. B89) so the style entry is re-used instead of duplicated.
$wsCost.Range("B88:C88").NumberFormat = """$""#,##0"
$wsCost.Range("B96").NumberFormat = """$""#,##0"

$wsCost.Range("B88").Formula = "=B54"
$wsCost.Range("C88").Formula = "=B55"
$wsCost.Range("B96").Formula = "=B87"
